# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-26 08:17:11
#
# This script applies the "Recorded By" re-ordering fixes, updates the
# Missing/Pending session roll-up counters, and flips three still-open
# sessions (one per group B2A/B2B/B2C) from "Pending" to "Not Recorded"
# (including re-applying the matching "Not Recorded" row formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Recorded By" (column G) text fix-ups: move "System" out of the
#    leading position of the comma separated list.
# ---------------------------------------------------------------------
$GChanges = @(
    @{Cell='G2'; Value='backup@backdoor.com, system, System'}
    @{Cell='G3'; Value='dnasr281@gmail.com, System'}
    @{Cell='G6'; Value='dnasr281@gmail.com, System'}
    @{Cell='G7'; Value='admin@admin.com, System'}
    @{Cell='G10'; Value='dnasr281@gmail.com, System'}
    @{Cell='G11'; Value='dnasr281@gmail.com, System'}
    @{Cell='G12'; Value='dnasr281@gmail.com, System'}
    @{Cell='G13'; Value='dnasr281@gmail.com, System'}
    @{Cell='G14'; Value='dnasr281@gmail.com, System'}
    @{Cell='G15'; Value='dnasr281@gmail.com, System'}
    @{Cell='G17'; Value='dnasr281@gmail.com, System'}
    @{Cell='G18'; Value='dnasr281@gmail.com, System'}
    @{Cell='G19'; Value='dnasr281@gmail.com, System'}
    @{Cell='G20'; Value='dnasr281@gmail.com, System'}
    @{Cell='G21'; Value='dnasr281@gmail.com, System'}
    @{Cell='G22'; Value='dnasr281@gmail.com, System'}
    @{Cell='G29'; Value='backup@backdoor.com, system, System'}
    @{Cell='G30'; Value='dnasr281@gmail.com, System'}
    @{Cell='G33'; Value='dnasr281@gmail.com, System'}
    @{Cell='G34'; Value='admin@admin.com, System'}
    @{Cell='G37'; Value='dnasr281@gmail.com, System'}
    @{Cell='G38'; Value='dnasr281@gmail.com, System'}
    @{Cell='G39'; Value='dnasr281@gmail.com, System'}
    @{Cell='G40'; Value='dnasr281@gmail.com, System'}
    @{Cell='G41'; Value='dnasr281@gmail.com, System'}
    @{Cell='G42'; Value='dnasr281@gmail.com, System'}
    @{Cell='G44'; Value='dnasr281@gmail.com, System'}
    @{Cell='G45'; Value='dnasr281@gmail.com, System'}
    @{Cell='G46'; Value='dnasr281@gmail.com, System'}
    @{Cell='G47'; Value='dnasr281@gmail.com, System'}
    @{Cell='G48'; Value='dnasr281@gmail.com, System'}
    @{Cell='G49'; Value='dnasr281@gmail.com, System'}
    @{Cell='G56'; Value='backup@backdoor.com, system, System'}
    @{Cell='G57'; Value='dnasr281@gmail.com, System'}
    @{Cell='G60'; Value='dnasr281@gmail.com, System'}
    @{Cell='G61'; Value='admin@admin.com, System'}
    @{Cell='G64'; Value='dnasr281@gmail.com, System'}
    @{Cell='G65'; Value='dnasr281@gmail.com, System'}
    @{Cell='G66'; Value='dnasr281@gmail.com, System'}
    @{Cell='G67'; Value='dnasr281@gmail.com, System'}
    @{Cell='G68'; Value='dnasr281@gmail.com, System'}
    @{Cell='G69'; Value='dnasr281@gmail.com, System'}
    @{Cell='G71'; Value='dnasr281@gmail.com, System'}
    @{Cell='G72'; Value='dnasr281@gmail.com, System'}
    @{Cell='G73'; Value='dnasr281@gmail.com, System'}
    @{Cell='G74'; Value='dnasr281@gmail.com, System'}
    @{Cell='G75'; Value='dnasr281@gmail.com, System'}
    @{Cell='G76'; Value='dnasr281@gmail.com, System'}
    @{Cell='G86'; Value='dnasr281@gmail.com, System'}
    @{Cell='G87'; Value='dnasr281@gmail.com, System'}
    @{Cell='G88'; Value='dnasr281@gmail.com, System'}
    @{Cell='G89'; Value='dnasr281@gmail.com, System'}
    @{Cell='G93'; Value='dnasr281@gmail.com, System'}
    @{Cell='G95'; Value='dnasr281@gmail.com, System'}
    @{Cell='G96'; Value='dnasr281@gmail.com, System'}
    @{Cell='G97'; Value='dnasr281@gmail.com, System'}
    @{Cell='G99'; Value='dnasr281@gmail.com, System'}
    @{Cell='G102'; Value='dnasr281@gmail.com, System'}
    @{Cell='G112'; Value='dnasr281@gmail.com, System'}
    @{Cell='G113'; Value='dnasr281@gmail.com, System'}
    @{Cell='G114'; Value='dnasr281@gmail.com, System'}
    @{Cell='G115'; Value='dnasr281@gmail.com, System'}
    @{Cell='G119'; Value='dnasr281@gmail.com, System'}
    @{Cell='G121'; Value='dnasr281@gmail.com, System'}
    @{Cell='G122'; Value='dnasr281@gmail.com, System'}
    @{Cell='G123'; Value='dnasr281@gmail.com, System'}
    @{Cell='G125'; Value='dnasr281@gmail.com, System'}
    @{Cell='G128'; Value='dnasr281@gmail.com, System'}
    @{Cell='G138'; Value='dnasr281@gmail.com, System'}
    @{Cell='G139'; Value='dnasr281@gmail.com, System'}
    @{Cell='G140'; Value='dnasr281@gmail.com, System'}
    @{Cell='G141'; Value='dnasr281@gmail.com, System'}
    @{Cell='G145'; Value='dnasr281@gmail.com, System'}
    @{Cell='G147'; Value='dnasr281@gmail.com, System'}
    @{Cell='G148'; Value='dnasr281@gmail.com, System'}
    @{Cell='G149'; Value='dnasr281@gmail.com, System'}
    @{Cell='G151'; Value='dnasr281@gmail.com, System'}
    @{Cell='G154'; Value='dnasr281@gmail.com, System'}
)

foreach ($change in $GChanges) {
    $ws.Range($change.Cell).Value = $change.Value
}

# ---------------------------------------------------------------------
# 2) Workbook level roll-up counters (Missing / Pending sessions).
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 6
$ws.Range("L8").Value = 30

# ---------------------------------------------------------------------
# 3) Per-group "Missing" / "Pending" counters (B2A, B2B, B2C) impacted
#    by the status flips below.
# ---------------------------------------------------------------------
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 4
$ws.Range("P16").Value = 2
$ws.Range("Q16").Value = 4
$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = 4

# ---------------------------------------------------------------------
# 4) Flip the still-open session rows 24 (B2A), 51 (B2B) and 78 (B2C)
#    from "Pending" to "Not Recorded", matching both the text and the
#    "Not Recorded" (pink) row styling used by rows 23, 50 and 77.
# ---------------------------------------------------------------------
$StatusFlipRows = @(
    @{Target=24; FormatSource=23},
    @{Target=51; FormatSource=50},
    @{Target=78; FormatSource=77}
)

foreach ($flip in $StatusFlipRows) {
    $target = $flip.Target
    $src = $flip.FormatSource

    # Update the status text first.
    $ws.Range("I$target").Value = "Not Recorded"

    # Re-apply the "Not Recorded" formatting (fill/font) from the row
    # above, which already carries the correct style, without touching
    # any of the underlying cell values.
    $srcRange = $ws.Range("A$src`:I$src")
    $dstRange = $ws.Range("A$target`:I$target")
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
